# Auto-generated edit script: updates pricing/profit columns (H-N)
# on the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to match the refreshed
# market-board snapshot pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 2951
$ws.Range("I16").Value = 2921.6
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 2921.6
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -2691.6
$ws.Range("N16").Value = -3460
$ws.Range("H18").Value = 29999.5
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H88").Value = 2500
$ws.Range("I88").Value = 2500
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 2500
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -2094
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 2500
$ws.Range("I91").Value = 2500
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 2500
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -1096
$ws.Range("N91").ClearContents()
$ws.Range("H106").Value = 2499.5
$ws.Range("I106").Value = 2499
$ws.Range("K106").Value = 2499
$ws.Range("M106").Value = -1868
$ws.Range("H107").Value = 3573.3635
$ws.Range("I107").Value = 3125.75
$ws.Range("K107").Value = 3125.75
$ws.Range("M107").Value = -1205.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 50226.145
$ws.Range("J31").Value = 111195.336
$ws.Range("L31").Value = 111195.336
$ws.Range("N31").Value = -111783.336
$ws.Range("H32").Value = 13895896
$ws.Range("I32").Value = 13895896
$ws.Range("K32").Value = 13895896
$ws.Range("M32").Value = -13895609
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("H61").Value = 26374374
$ws.Range("I61").Value = 83334650
$ws.Range("K61").Value = 83334650
$ws.Range("M61").Value = -83334438
$ws.Range("H74").Value = 7818667
$ws.Range("I74").Value = 12500823
$ws.Range("K74").Value = 12500823
$ws.Range("M74").Value = -12499949
$ws.Range("H77").Value = 7818667
$ws.Range("I77").Value = 12500823
$ws.Range("K77").Value = 62504115
$ws.Range("M77").Value = -62499747
$ws.Range("H132").Value = 8282.056
$ws.Range("I132").Value = 4339.3076
$ws.Range("K132").Value = 13017.9228
$ws.Range("M132").Value = -10487.9228
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 26374374
$ws.Range("I136").Value = 83334650
$ws.Range("K136").Value = 250003950
$ws.Range("M136").Value = -250001400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3464.3333
$ws.Range("I86").Value = 3214.4443
$ws.Range("J86").Value = 3839.1667
$ws.Range("K86").Value = 3214.4443
$ws.Range("L86").Value = 3839.1667
$ws.Range("M86").Value = -2091.4443
$ws.Range("N86").Value = -6085.1667
$ws.Range("H89").Value = 3464.3333
$ws.Range("I89").Value = 3214.4443
$ws.Range("J89").Value = 3839.1667
$ws.Range("K89").Value = 16072.2215
$ws.Range("L89").Value = 19195.8335
$ws.Range("M89").Value = -10456.2215
$ws.Range("N89").Value = -30427.8335
$ws.Range("H94").Value = 1455.6154
$ws.Range("I94").Value = 1410.25
$ws.Range("K94").Value = 1410.25
$ws.Range("M94").Value = -959.25
$ws.Range("H96").Value = 47158
$ws.Range("I96").Value = 17499.75
$ws.Range("J96").Value = 70884.60000000001
$ws.Range("K96").Value = 17499.75
$ws.Range("L96").Value = 70884.60000000001
$ws.Range("M96").Value = -14753.75
$ws.Range("N96").Value = -76376.60000000001
$ws.Range("H134").Value = 40035.883
$ws.Range("I134").Value = 1636.76
$ws.Range("K134").Value = 4910.28
$ws.Range("M134").Value = -2375.28

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1016072.06
$ws.Range("J31").Value = 2335015
$ws.Range("L31").Value = 2335015
$ws.Range("N31").Value = -2335605
$ws.Range("H34").Value = 1016072.06
$ws.Range("J34").Value = 2335015
$ws.Range("L34").Value = 2335015
$ws.Range("N34").Value = -2335419
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H51").Value = 45583.168
$ws.Range("I51").Value = 20142.572
$ws.Range("J51").Value = 81200
$ws.Range("K51").Value = 20142.572
$ws.Range("L51").Value = 81200
$ws.Range("M51").Value = -19406.572
$ws.Range("N51").Value = -82672
$ws.Range("H61").Value = 45583.168
$ws.Range("I61").Value = 20142.572
$ws.Range("J61").Value = 81200
$ws.Range("K61").Value = 20142.572
$ws.Range("L61").Value = 81200
$ws.Range("M61").Value = -19794.572
$ws.Range("N61").Value = -81896
$ws.Range("I62").Value = 2165.6667
$ws.Range("J62").Value = 2232.3333
$ws.Range("K62").Value = 2165.6667
$ws.Range("L62").Value = 2232.3333
$ws.Range("M62").Value = -1541.6667
$ws.Range("N62").Value = -3480.3333
$ws.Range("H64").Value = 108000
$ws.Range("J64").Value = 108000
$ws.Range("L64").Value = 108000
$ws.Range("N64").Value = -108496
$ws.Range("I65").Value = 2165.6667
$ws.Range("J65").Value = 2232.3333
$ws.Range("K65").Value = 10828.3335
$ws.Range("L65").Value = 11161.6665
$ws.Range("M65").Value = -7708.333500000001
$ws.Range("N65").Value = -17401.6665
$ws.Range("H67").Value = 108000
$ws.Range("J67").Value = 108000
$ws.Range("L67").Value = 108000
$ws.Range("N67").Value = -109716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2024223.6
$ws.Range("J81").Value = 2529154.5
$ws.Range("L81").Value = 7587463.5
$ws.Range("N81").Value = -7589709.5
$ws.Range("H84").Value = 2024223.6
$ws.Range("J84").Value = 2529154.5
$ws.Range("L84").Value = 22762390.5
$ws.Range("N84").Value = -22773622.5
$ws.Range("H86").Value = 122.166664
$ws.Range("I86").Value = 123
$ws.Range("J86").Value = 120.5
$ws.Range("K86").Value = 369
$ws.Range("L86").Value = 361.5
$ws.Range("M86").Value = 817
$ws.Range("N86").Value = -2733.5
$ws.Range("H89").Value = 122.166664
$ws.Range("I89").Value = 123
$ws.Range("J89").Value = 120.5
$ws.Range("K89").Value = 1107
$ws.Range("L89").Value = 1084.5
$ws.Range("M89").Value = 4821
$ws.Range("N89").Value = -12940.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 500000
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H70").Value = 3333
$ws.Range("I70").Value = 3333
$ws.Range("K70").Value = 3333
$ws.Range("M70").Value = -3063
$ws.Range("H73").Value = 3333
$ws.Range("I73").Value = 3333
$ws.Range("K73").Value = 3333
$ws.Range("M73").Value = -2397
$ws.Range("H97").Value = 1052.0526
$ws.Range("J97").Value = 2548
$ws.Range("L97").Value = 2548
$ws.Range("N97").Value = -3540
$ws.Range("H139").Value = 88000
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 5020.5
$ws.Range("I45").Value = 5020.5
$ws.Range("K45").Value = 5020.5
$ws.Range("M45").Value = -4613.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 62190
$ws.Range("J92").Value = 62190
$ws.Range("L92").Value = 62190
$ws.Range("N92").Value = -67182
$ws.Range("H95").Value = 98311.336
$ws.Range("J95").Value = 98311.336
$ws.Range("L95").Value = 98311.336
$ws.Range("N95").Value = -103803.336
$ws.Range("H132").Value = 3069.9
$ws.Range("I132").Value = 1139.8
$ws.Range("K132").Value = 3419.4
$ws.Range("M132").Value = -889.3999999999996
$ws.Range("H136").Value = 6238.125
$ws.Range("I136").Value = 6676.905
$ws.Range("J136").Value = 3166.6667
$ws.Range("K136").Value = 20030.715
$ws.Range("L136").Value = 9500.000100000001
$ws.Range("M136").Value = -17480.715
$ws.Range("N136").Value = -14600.0001

